$wb = $excel.ActiveWorkbook

# --- Sheet 2: "工作表 1 - 题目一_一般边界值法" ---
# Fill in rows 27-29, columns A-D with boundary-value test data.
$wsBoundary = $wb.Worksheets.Item(2)

$wsBoundary.Range("A27").Value = 1900
$wsBoundary.Range("B27").Value = 12
$wsBoundary.Range("C27").Value = 31
$wsBoundary.Range("D27").Value = 1

$wsBoundary.Range("A28").Value = 2000
$wsBoundary.Range("B28").Value = 12
$wsBoundary.Range("C28").Value = 31
$wsBoundary.Range("D28").Value = 7

$wsBoundary.Range("A29").Value = 2100
$wsBoundary.Range("B29").Value = 12
$wsBoundary.Range("C29").Value = 31
$wsBoundary.Range("D29").Value = 5

# --- Sheet 3: "工作表 2 - 题目一_弱健壮等价类法" ---
# Clear the contents (keep formatting) of rows 42-46, columns A-F.
$wsRobust = $wb.Worksheets.Item(3)
$wsRobust.Range("A42:F46").ClearContents()

# --- View state: make sheet 3 the active / selected tab ---
$wsBoundary.Activate()
$wsBoundary.Range("A27:D29").Select()

$wb.Windows.Item(1).ScrollRow = 17
$wb.Windows.Item(1).ScrollColumn = 2

$wsRobust.Activate()
$wsRobust.Range("A42:F46").Select()

$wb.Windows.Item(1).ScrollRow = 35
$wb.Windows.Item(1).ScrollColumn = 2
